# This script rotates the data of rows 2, 3 and 4 (row 3 -> row 2, row 4 -> row 3,
# row 2 -> row 4) while leaving the columns that are identical across the three
# rows (C, J, K, L, N, O, T, U, V, W, X, Y, AD, AE, AF, AG, AH-AS, AT, AU, AV, AW,
# AX, AY, ...) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: becomes the former row 3 data ----
$ws.Range("A2").Value = 130826010
$ws.Range("B2").Value = 91808
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = "Ullticka"
$ws.Range("G2").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H2").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M2").ClearContents()
$ws.Range("P2").Value = "Flinktorpet, Flinktorpet, Jmt"
$ws.Range("Q2").Value = 460971
$ws.Range("R2").Value = 7039688
$ws.Range("Z2").Value = "10:47"
$ws.Range("AB2").Value = "10:47"
$ws.Range("AC2").ClearContents()

# ---- Row 3: becomes the former row 4 data ----
$ws.Range("A3").Value = 130826784
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("M3").Value = "färska spår"
$ws.Range("P3").Value = "Brännan, Kälom, Offerdal, Jmt"
$ws.Range("Q3").Value = 461233
$ws.Range("R3").Value = 7039438
$ws.Range("Z3").Value = "11:37"
$ws.Range("AB3").Value = "11:37"
$ws.Range("AC3").Value = "Födosök barkfläk"

# ---- Row 4: becomes the former row 2 data ----
$ws.Range("A4").Value = 130825823
$ws.Range("B4").Value = 57881
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "äldre spår"
$ws.Range("P4").Value = "Flinktorpet, Kälom, Offerdal, Jmt"
$ws.Range("Q4").Value = 460947
$ws.Range("R4").Value = 7039711
$ws.Range("Z4").Value = "10:38"
$ws.Range("AB4").Value = "10:38"
$ws.Range("AC4").Value = "Födosökshål på äldre döende gran."
